$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.742.92'
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').Value = '1.817.04'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.569'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.60%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '35.06'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.300'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '2.078.39'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.08%  '
$ws.Range('D14').Value = '1.821.19'
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '34.721.49'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('E17').Value = '  +2.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.30%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '171.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.03%  '
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +4.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.118'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('D36').Value = '1.422.07'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.683'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.71%  '
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.958'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('E45').Value = '  +3.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0518'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').Value = '1.978.30'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  +1.91%  '
